$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H, shifting old H:J to I:K
$ws.Range("H1").EntireColumn.Insert()

# Match column G's (bestFit) width on the freshly inserted column H
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# New header cell for the "Instructor" column
$ws.Range("H1").Value = "Instructor"

# Template cell already carrying the "Arial 11 black" style used
# throughout the data rows (column G uses it for many rows) - copy
# just the formatting onto every new "OJ" cell below the header.
$styleTemplate = $ws.Cells.Item(2, 7)
$styleTemplate.Copy()

for ($r = 2; $r -le 45; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Value = "OJ"
    $cell.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# Restore the cursor/selection to where the author left it
[void]$ws.Range("G20").Select()
